# Append four new rows of hourly charging data (two days x two stations)
# to Sheet1, mirroring the existing table layout (columns A:Z).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 120; Date = 45960; Label = "四方坪站充电量(kw)"; Values = @(989.6149999999999, 734.03, 445.57700000000006, 311, 409.13800000000003, 608.76799999999992, 280.58100000000002, 115.31300000000002, 181.61100000000002, 168.91200000000001, 135.65200000000002, 124.774, 806.14900000000011, 1129.0839999999998, 556.45500000000004, 414.06400000000008, 335.22399999999999, 133.941, 42.72, 282.62299999999999, 86.783000000000001, 213.98999999999998, 23.824000000000002, 33.78) }
    @{ Row = 121; Date = 45960; Label = "高岭站充电量(kw)"; Values = @(293.13900000000001, 285.20399999999995, 118.218, 60.241999999999997, 33.692, 185.57900000000004, 159.172, 132.84700000000001, 141.11199999999999, 278.24400000000003, 150.92099999999999, 188.73699999999997, 387.25900000000007, 233.34200000000001, 313.71899999999999, 177.96899999999999, 425.65800000000002, 128.017, 36.773000000000003, 190.48299999999998, 52.43, 78.98, 15.257, 42.008000000000003) }
    @{ Row = 122; Date = 45961; Label = "四方坪站充电量(kw)"; Values = @(616.39100000000008, 1017.5880000000001, 397.84200000000004, 651.91200000000003, 129.25, 793.01699999999994, 354.01799999999997, 170.83500000000001, 119.13, 154.23099999999999, 163.32999999999998, 179.72000000000003, 715.08199999999999, 1469.6569999999997, 598.197, 548.22800000000007, 451.88099999999997, 364.17399999999998, 87.78, 43.012999999999998, 265.32400000000001, 46.18, 71.62, 132.85) }
    @{ Row = 123; Date = 45961; Label = "高岭站充电量(kw)"; Values = @(502.28399999999999, 317.34600000000006, 145.398, 87.664000000000001, 60.992999999999995, 69.507999999999996, 241.77199999999999, 189.33199999999999, 216.34899999999999, 182.83499999999998, 181.51, 176.26400000000001, 306.68099999999998, 632.60599999999988, 683.46100000000001, 220.37299999999996, 334.89300000000003, 189.06200000000001, 49.533999999999999, 18.952999999999999, 87.822000000000003, 26.047000000000001, 70.031000000000006, 13.016) }
)

foreach ($entry in $newRows) {
    $r = $entry.Row

    $ws.Cells.Item($r, 1).Value = $entry.Date
    $ws.Cells.Item($r, 2).Value = $entry.Label

    $col = 3
    foreach ($v in $entry.Values) {
        $ws.Cells.Item($r, $col).Value = $v
        $col = $col + 1
    }
}

# Match the author's final viewport/selection state: scrolled down so row 89
# is at the top, with E125 as the active (selected) cell.
$ws.Range("A89").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 89
$ws.Range("E125").Select() | Out-Null
